$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I: "category" (normal), pushing the old I (date),
# J (legislator_name) and K (legislator_id) columns one to the right.
$ws.Columns.Item(9).Insert()

# Insert two more columns after L (legislator_id) for source_file / index,
# using Insert() (rather than writing into blank cells) so the new cells
# pick up the same header/data styles (s=1 / s=2) as their neighbours.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()

# Header row
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data rows 2-11: new "category" column = "normal"
$ws.Range("I2:I11").Value = "normal"

# New trailing columns: source_file = "tmpf41", index = same as column A
$ws.Range("M2:M11").Value = "tmpf41"

for ($r = 2; $r -le 11; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 14).Value = $idx
}
